$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-01-20 01:27:18"
}
